# Update countries & provincias Spain
# Applies updated COVID-19 case figures for several countries. A few
# countries' totals overtook the country immediately above them in the
# (descending, by total cases) list, so their row now carries the new
# figures while the displaced country keeps its previous figures one
# row further down (mirroring the source data re-sort).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Belgica (row 13) - direct refresh of all figures, still correctly ranked
$ws.Range("B13").Value = 34809
$ws.Range("C13").Value = 1236
$ws.Range("D13").Value = 7562
$ws.Range("E13").Value = 22390
$ws.Range("F13").Value = 1182
$ws.Range("G13").Value = 417
$ws.Range("H13").Value = 4857

# Austria (row 20) - only totals/new-cases/active refreshed
$ws.Range("B20").Value = 14412
$ws.Range("C20").Value = 62
$ws.Range("E20").Value = 5921

# Indonesia overtakes Emiratos Arabes Unidos -> row 40 becomes Indonesia
# (with fresh data) and Emiratos Arabes Unidos shifts to row 41 keeping
# its previous data.
$ws.Range("A40").Value = "Indonesia"
$ws.Range("B40").Value = 5516
$ws.Range("C40").Value = 380
$ws.Range("D40").Value = 548
$ws.Range("E40").Value = 4470
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 29
$ws.Range("H40").Value = 498

$ws.Range("A41").Value = "Emiratos Arabes Unidos"
$ws.Range("B41").Value = 5365
$ws.Range("C41").Value = 0
$ws.Range("D41").Value = 1034
$ws.Range("E41").Value = 4298
$ws.Range("F41").Value = 1
$ws.Range("G41").Value = 0
$ws.Range("H41").Value = 33

# Moldavia (row 59) - only active cases/recuperados refreshed
$ws.Range("D59").Value = 235
$ws.Range("E59").Value = 1763

# Estonia overtakes Irak -> row 66 becomes Estonia (fresh data), and
# Irak, Kuwait, Nueva Zelanda each shift down one row keeping their
# previous data.
$ws.Range("A66").Value = "Estonia"
$ws.Range("B66").Value = 1434
$ws.Range("C66").Value = 34
$ws.Range("D66").Value = 133
$ws.Range("E66").Value = 1265
$ws.Range("F66").Value = 10
$ws.Range("G66").Value = 1
$ws.Range("H66").Value = 36

$ws.Range("A67").Value = "Irak"
$ws.Range("B67").Value = 1415
$ws.Range("C67").Value = 0
$ws.Range("D67").Value = 812
$ws.Range("E67").Value = 524
$ws.Range("F67").Value = 0
$ws.Range("G67").Value = 0
$ws.Range("H67").Value = 79

$ws.Range("A68").Value = "Kuwait"
$ws.Range("B68").Value = 1405
$ws.Range("C68").Value = 0
$ws.Range("D68").Value = 206
$ws.Range("E68").Value = 1196
$ws.Range("F68").Value = 31
$ws.Range("G68").Value = 0
$ws.Range("H68").Value = 3

$ws.Range("A69").Value = "Nueva Zelanda"
$ws.Range("B69").Value = 1401
$ws.Range("C69").Value = 15
$ws.Range("D69").Value = 770
$ws.Range("E69").Value = 622
$ws.Range("F69").Value = 3
$ws.Range("G69").Value = 0
$ws.Range("H69").Value = 9

# Isla de Man overtakes Vietnam -> row 114 becomes Isla de Man (fresh
# data), and Vietnam, Consejo Danes para los Refugiados each shift down
# one row keeping their previous data.
$ws.Range("A114").Value = "Isla de Man"
$ws.Range("B114").Value = 283
$ws.Range("C114").Value = 27
$ws.Range("D114").Value = 153
$ws.Range("E114").Value = 126
$ws.Range("F114").Value = 13
$ws.Range("G114").Value = 0
$ws.Range("H114").Value = 4

$ws.Range("A115").Value = "Vietnam"
$ws.Range("B115").Value = 268
$ws.Range("C115").Value = 0
$ws.Range("D115").Value = 171
$ws.Range("E115").Value = 97
$ws.Range("F115").Value = 8
$ws.Range("G115").Value = 0
$ws.Range("H115").Value = 0

$ws.Range("A116").Value = "Consejo Danes para los Refugiados"
$ws.Range("B116").Value = 267
$ws.Range("C116").Value = 13
$ws.Range("D116").Value = 23
$ws.Range("E116").Value = 222
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 1
$ws.Range("H116").Value = 22
